# Update ListShip sequence diagram (slide 1 / sldId 267):
# reposition/resize the three lifeline bars and the connectors attached
# to the "list" interaction so the arrows line up with their new extents.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Rectangle 224 (id 225) - lifeline segment -----------------------------
$shp = $s.Shapes.Item("Rectangle 224")
$shp.Top    = 530.0504150390625
$shp.Height = 296.0645751953125

# --- Rectangle 225 (id 226) - lifeline segment -----------------------------
$shp = $s.Shapes.Item("Rectangle 225")
$shp.Left   = 183.50662231445312
$shp.Top    = 868.2361450195312
$shp.Width  = 14.970630645751953
$shp.Height = 271.3208923339844

# --- Rectangle 226 (id 227) - lifeline segment -----------------------------
$shp = $s.Shapes.Item("Rectangle 226")
$shp.Left   = 183.50662231445312
$shp.Top    = 1219.477294921875
$shp.Width  = 13.587165832519531
$shp.Height = 272.6754455566406

# --- Straight Arrow Connector 229 (id 230) ---------------------------------
$shp = $s.Shapes.Item("Straight Arrow Connector 229")
$shp.Left   = 194.92771911621094
$shp.Top    = 871.7572021484375
$shp.Width  = 201.73536682128906
$shp.Height = 0

# --- Straight Arrow Connector 291 (id 292) ---------------------------------
$shp = $s.Shapes.Item("Straight Arrow Connector 291")
$shp.VerticalFlip = -1
$shp.Top    = 530.1129150390625
$shp.Width  = 104.5572509765625
$shp.Height = 0.9022047519683838

# --- Straight Arrow Connector 187 (id 188) ---------------------------------
$shp = $s.Shapes.Item("Straight Arrow Connector 187")
$shp.Width  = 90.22535705566406
